$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 - new data (A, B, E, F are literal values; C, D, G, H are formulas)
$ws.Range("A8").Value = 1440
$ws.Range("B8").Value = 1440
$ws.Range("C8").Formula = "=B8/3"
$ws.Range("D8").Formula = "=2*B8/3"
$ws.Range("E8").Value = 620
$ws.Range("F8").Value = 620
$ws.Range("G8").Formula = "=B8/2-F8/2"
$ws.Range("H8").Formula = "=A8/2-E8/2"
$ws.Range("H8").Style = "Normal"

# Row 9 - new data (A, B, E, F are literal values; C, D, G, H are formulas)
$ws.Range("A9").Value = 540
$ws.Range("B9").Value = 960
$ws.Range("C9").Formula = "=B9/3"
$ws.Range("D9").Formula = "=2*B9/3"
$ws.Range("E9").Value = 420
$ws.Range("F9").Value = 420
$ws.Range("G9").Formula = "=B9/2-F9/2"
$ws.Range("H9").Formula = "=A9/2-E9/2"
$ws.Range("H9").Style = "Normal"

# Update the selection to match the final state
$ws.Range("C9").Select()
